$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# "Generate Report for Handback" - the localization-status report picked
# up a new handback for a2a8cc29-6ddb-4884-bd6c-6f5f0babf4a9 on both the
# zh-cn and de-de target sheets (row 7): the "Latest Target File" /
# "Latest Handback DateTime" / "Error Detail" columns get filled in, and
# a hyperlink is added on the "Latest Target File" cell (column I).
# -----------------------------------------------------------------------

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c12bc9baeca541c485b084d99484a98f09ca7c00/e2e/a2a8cc29-6ddb-4884-bd6c-6f5f0babf4a9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0de719d63748867d5270ac7a4fce1803a318572/e2e/a2a8cc29-6ddb-4884-bd6c-6f5f0babf4a9.md."
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0de719d63748867d5270ac7a4fce1803a318572/e2e/a2a8cc29-6ddb-4884-bd6c-6f5f0babf4a9.md"
$displayName = "a2a8cc29-6ddb-4884-bd6c-6f5f0babf4a9.md"

$sheetNames = @("zh-cn", "de-de")
$targetFileNames = @(
    "a2a8cc29-6ddb-4884-bd6c-6f5f0babf4a9.e6cc5951783652542b007f57025be0d5093f5559.zh-cn.xlf",
    "a2a8cc29-6ddb-4884-bd6c-6f5f0babf4a9.e6cc5951783652542b007f57025be0d5093f5559.de-de.xlf"
)
$handbackDateTimes = @("2016-09-04 02:59:23", "2016-09-04 02:59:30")

for ($i = 0; $i -lt $sheetNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])

    # I7: "Latest Target File" -> becomes a hyperlink to the handed-back
    # source markdown file (same file/URL already referenced by A7).
    $i7 = $ws.Range("I7")
    $i7.Value = $displayName
    $ws.Hyperlinks.Add($i7, $targetUrl, "", "", $displayName) | Out-Null
    $i7.Font.Underline = 2
    $i7.Font.Color = 15570276

    # J7: "Latest Handback File" -> the actual xliff file name handed back.
    $ws.Range("J7").Value = $targetFileNames[$i]

    # K7: "Latest Handback DateTime" -> timestamp of the handback.
    $ws.Range("K7").Value = $handbackDateTimes[$i]

    # P7: "Error Detail" -> version mismatch warning.
    $ws.Range("P7").Value = $errorDetail
}
